# Week 3 Progress Report update:
# Add a new "Work Done Upto" entry for Atharva Kulkarni (row 4, column D)
# describing the work completed this week.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = "Studied basics of UI Design for Flutter"

# Move the active selection to the newly filled cell.
$ws.Range("D4").Select() | Out-Null
